# feature: implement create code item
#
# Fixes the item_code/category of the first product row and appends three
# new "Alimentos y bebidas" stock items (with their own generated item
# codes) to the ProductData sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductData")

# Reuse the creation_date column's existing number format for the new rows.
$dateFmt = $ws.Cells.Item(2, 7).NumberFormat()

# Row 2: correct the item code scheme and the "Alimento y bebidas" typo,
# and refresh the creation timestamp.
$ws.Cells.Item(2, 1).Value = "1AYB"
$ws.Cells.Item(2, 2).Value = "Alimentos y bebidas"
$ws.Cells.Item(2, 7).Value = 45798.94746832781

# Newly created stock items.
$newItems = @(
    @("2AYB", "Alimentos y bebidas", "Aceite de oliva x 1 Litro", 3, 8000, 12000, 45798.9474683279),
    @("3AYB", "Alimentos y bebidas", "Leche Colanta x 1 litro", 8, 2500, 3500, 45798.94746832792),
    @("4AYB", "Alimentos y bebidas", "Pan Integral Unidad", 4, 1000, 1500, 45798.94746832793)
)

$r = 3
foreach ($item in $newItems) {
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
    $ws.Cells.Item($r, 6).Value = $item[5]
    $ws.Cells.Item($r, 7).Value = $item[6]
    $ws.Cells.Item($r, 7).NumberFormat = $dateFmt
    $r = $r + 1
}
